$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 1.164924666666667
$ws.Range("H2").Value = 3.494774
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.017765
$ws.Range("N2").Value = 0.053295
$ws.Range("O2").Value = 0.005225147533577419
$ws.Range("P2").Value = 0.005225147533577419
$ws.Range("Q2").Value = 0.02069488670333333
$ws.Range("R2").Value = 0.18625398033
$ws.Range("S2").Value = 0.005225147533577419
$ws.Range("T2").Value = 0.005225147533577419

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 1.164924666666667
$ws.Range("H3").Value = 3.494774
$ws.Range("M3").Value = 0.8788360000000001
$ws.Range("N3").Value = 2.636508
$ws.Range("O3").Value = 0.258488474968705
$ws.Range("P3").Value = 0.258488474968705
$ws.Range("Q3").Value = 1.023777734354667
$ws.Range("R3").Value = 9.213999609192001
$ws.Range("S3").Value = 0.258488474968705
$ws.Range("T3").Value = 0.258488474968705

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 1.164924666666667
$ws.Range("H4").Value = 3.494774
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 2.503303
$ws.Range("N4").Value = 7.509909
$ws.Range("O4").Value = 0.7362863774977175
$ws.Range("P4").Value = 0.7362863774977175
$ws.Range("Q4").Value = 2.916159412840666
$ws.Range("R4").Value = 26.245434715566
$ws.Range("S4").Value = 0.7362863774977175
$ws.Range("T4").Value = 0.7362863774977175
